# Add "IP Mobility on EPG" locale strings to the "locale" worksheet.
# New rows 169-175 are appended (Key / en / ko-KR / ko columns A-D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 169
$ws.Cells.Item(169, 1).Value = "IP Mobility"
$ws.Cells.Item(169, 2).Value = "IP Mobility"
$ws.Cells.Item(169, 3).Value = "IP 모빌리티"
$ws.Cells.Item(169, 4).Value = "IP 모빌리티"

# row 170
$ws.Cells.Item(170, 1).Value = "Register IPs"
$ws.Cells.Item(170, 2).Value = "Register IPs"
$ws.Cells.Item(170, 3).Value = "IP(구간) 등록"
$ws.Cells.Item(170, 4).Value = "IP(구간) 등록"

# row 171
$ws.Cells.Item(171, 1).Value = "Mobility IPs"
$ws.Cells.Item(171, 2).Value = "Mobility IPs"
$ws.Cells.Item(171, 3).Value = "IP 모빌리티 리스트"
$ws.Cells.Item(171, 4).Value = "IP 모빌리티 리스트"

# row 172
$ws.Cells.Item(172, 1).Value = "IP Start"
$ws.Cells.Item(172, 2).Value = "IP Start"
$ws.Cells.Item(172, 3).Value = "시작 IP"
$ws.Cells.Item(172, 4).Value = "시작 IP"

# row 173
$ws.Cells.Item(173, 1).Value = "IP End"
$ws.Cells.Item(173, 2).Value = "IP End"
$ws.Cells.Item(173, 3).Value = "종료 IP"
$ws.Cells.Item(173, 4).Value = "종료 IP"

# row 174
$ws.Cells.Item(174, 1).Value = "Description"
$ws.Cells.Item(174, 2).Value = "Description"
$ws.Cells.Item(174, 3).Value = "부가설명"
$ws.Cells.Item(174, 4).Value = "부가설명"

# row 175
$ws.Cells.Item(175, 1).Value = "Use EPG Subnet"
$ws.Cells.Item(175, 2).Value = "Use EPG Subnet"
$ws.Cells.Item(175, 3).Value = "EPG 서브넷 사용"
$ws.Cells.Item(175, 4).Value = "EPG 서브넷 사용"

# Update the view so the newly added row is the active selection, matching
# the author's saved sheet view state after entering the data.
$excel.ActiveWindow.ScrollRow = 151
$ws.Range("C176").Select()
